$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-4) get cyclically rotated for columns D, J, K, L, M, O, P:
#   new row2 = old row3, new row3 = old row4, new row4 = old row2
# Capture old values first (using Value2, which reliably returns primitives)
$r2D = $ws.Range("D2").Value2
$r2J = $ws.Range("J2").Value2
$r2K = $ws.Range("K2").Value2
$r2L = $ws.Range("L2").Value2
$r2M = $ws.Range("M2").Value2
$r2O = $ws.Range("O2").Value2
$r2P = $ws.Range("P2").Value2

$r3D = $ws.Range("D3").Value2
$r3J = $ws.Range("J3").Value2
$r3K = $ws.Range("K3").Value2
$r3L = $ws.Range("L3").Value2
$r3M = $ws.Range("M3").Value2
$r3O = $ws.Range("O3").Value2
$r3P = $ws.Range("P3").Value2

$r4D = $ws.Range("D4").Value2
$r4J = $ws.Range("J4").Value2
$r4K = $ws.Range("K4").Value2
$r4L = $ws.Range("L4").Value2
$r4M = $ws.Range("M4").Value2
$r4O = $ws.Range("O4").Value2
$r4P = $ws.Range("P4").Value2

# Write new row 2 (= old row 3)
$ws.Range("D2").Value2 = $r3D
$ws.Range("J2").Value2 = $r3J
$ws.Range("K2").Value2 = $r3K
$ws.Range("L2").Value2 = $r3L
$ws.Range("M2").Value2 = $r3M
$ws.Range("O2").Value2 = $r3O
$ws.Range("P2").Value2 = $r3P

# Write new row 3 (= old row 4)
$ws.Range("D3").Value2 = $r4D
$ws.Range("J3").Value2 = $r4J
$ws.Range("K3").Value2 = $r4K
$ws.Range("L3").Value2 = $r4L
$ws.Range("M3").Value2 = $r4M
$ws.Range("O3").Value2 = $r4O
$ws.Range("P3").Value2 = $r4P

# Write new row 4 (= old row 2)
$ws.Range("D4").Value2 = $r2D
$ws.Range("J4").Value2 = $r2J
$ws.Range("K4").Value2 = $r2K
$ws.Range("L4").Value2 = $r2L
$ws.Range("M4").Value2 = $r2M
$ws.Range("O4").Value2 = $r2O
$ws.Range("P4").Value2 = $r2P
